$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the "in progress" notes back to their plain data values
$ws.Range("A1").Value = "Donnée A1"
$ws.Range("B2").Value = "Donnée B2"
$ws.Range("C3").Value = "Donnée C3"

# Row 5 had been accidentally filled with row 1's data; restore its own values
$ws.Range("B5").Value = "Donnée B5"
$ws.Range("C5").Value = "Donnée C5"

# Move the active selection to C5
$ws.Range("C5").Select()
